# Apply the 'EventRegistrationForAMember' feature changes:
#  - Insert a new worksheet "createRegionalEvent" right before "sqlCount"
#    with header "eventName" / value "TestAutomation20191204134548".
#  - Make "addBrandNewVisitor" the active/selected tab (was
#    "addVisitorForExistingIndividual").

$wb = $excel.ActiveWorkbook

$sqlCountSheet = $wb.Worksheets.Item("sqlCount")

$newSheet = $wb.Worksheets.Add($sqlCountSheet)
$newSheet.Name = "createRegionalEvent"
$newSheet.Range("A1").Value = "eventName"
$newSheet.Range("A2").Value = "TestAutomation20191204134548"

$wb.Worksheets.Item("addBrandNewVisitor").Activate()
